# Update "想去人数" (attendance count) figures on the "展览" and "全部类型"
# sheets to match the refreshed data pulled from the source site.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2895
    $ws.Range("F4").Value = 100
    $ws.Range("F5").Value = 6700
    $ws.Range("F6").Value = 1626
}
